$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "1.002", "22.216.44") but must stay
# plain text like the source data, so force text format before assigning,
# then restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.216.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.557.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3805'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3306'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.49'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.142'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07392'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.842'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.756'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.552.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001076'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '86.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06650'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.405'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.217.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.277'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.564'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.938'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.728.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.090'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.917'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.907'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.362'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08216'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06320'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.16%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02339'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.322'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  -5.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.236'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6064'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.09%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.745'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5879'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.60'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.968'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.178'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07063'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.86%  '
